$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("G9").Value = 1.75
$ws.Range("H9").Value = 3.05
$ws.Range("I9").Value = 5.3
$ws.Range("L9").Value = 1.39
$ws.Range("M9").Value = 2.57
$ws.Range("N9").Value = 2.12
$ws.Range("O9").Value = 1.57
$ws.Range("P9").Value = 1.45
$ws.Range("Q9").Value = 2.4
$ws.Range("R9").Value = 1.93
$ws.Range("S9").Value = 1.7
$ws.Range("U9").Value = 7.3
$ws.Range("V9").Value = 8.25
$ws.Range("W9").Value = 14
$ws.Range("X9").Value = 16
$ws.Range("Y9").Value = 32
$ws.Range("Z9").Value = 7
$ws.Range("AA9").Value = 6.1
$ws.Range("AB9").Value = 16.5
$ws.Range("AC9").Value = 100
$ws.Range("AD9").Value = 900
$ws.Range("AE9").Value = 12
$ws.Range("AF9").Value = 32
$ws.Range("AG9").Value = 17
$ws.Range("AH9").Value = 120
$ws.Range("AI9").Value = 65

# Row 11
$ws.Range("G11").Value = 1.62
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 5.5
$ws.Range("AE11").Value = 15
$ws.Range("AF11").Value = 29

# Row 12
$ws.Range("I12").Value = 6.25
$ws.Range("N12").Value = 1.7
$ws.Range("O12").Value = 2.1
$ws.Range("U12").Value = 7
$ws.Range("W12").Value = 10
$ws.Range("AB12").Value = 19

# Row 23
$ws.Range("G23").Value = 2.55
$ws.Range("H23").Value = 3.25
$ws.Range("I23").Value = 2.8
$ws.Range("J23").Value = 1.07
$ws.Range("K23").Value = 9
$ws.Range("L23").Value = 1.33
$ws.Range("M23").Value = 3.25
$ws.Range("P23").Value = 1.4
$ws.Range("Q23").Value = 2.75
$ws.Range("AE23").Value = 8.5
$ws.Range("AF23").Value = 13

# Row 27
$ws.Range("G27").Value = 2
$ws.Range("I27").Value = 3.25
$ws.Range("K27").Value = 15
$ws.Range("R27").Value = 1.53
$ws.Range("S27").Value = 2.38
$ws.Range("AE27").Value = 15

# Row 29
$ws.Range("G29").Value = 2.65
$ws.Range("H29").Value = 2.75
$ws.Range("I29").Value = 2.87
$ws.Range("T29").Value = 6.3
$ws.Range("U29").Value = 11.75
$ws.Range("V29").Value = 10.5
$ws.Range("W29").Value = 32
$ws.Range("X29").Value = 28
$ws.Range("Z29").Value = 6.1
$ws.Range("AC29").Value = 100
$ws.Range("AE29").Value = 7
$ws.Range("AF29").Value = 13.5
$ws.Range("AG29").Value = 10.75
$ws.Range("AH29").Value = 37
$ws.Range("AI29").Value = 30
$ws.Range("AJ29").Value = 45
